# Auto-generated edit script applying numeric corrections to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H:N)
# across several job sheets, per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 1047.4
$ws.Range("I11").Value = 1047.4
$ws.Range("K11").Value = 1047.4
$ws.Range("M11").Value = -907.4000000000001
# Row 17
$ws.Range("H17").Value = 232.20512
$ws.Range("J17").Value = 232.20512
$ws.Range("L17").Value = 696.61536
$ws.Range("N17").Value = -1032.61536
# Row 32
$ws.Range("H32").Value = 1967.3334
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1967.3334
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1967.3334
$ws.Range("N32").Value = -2619.3334
$ws.Range("M32").ClearContents()
# Row 112
$ws.Range("H112").Value = 2201.25
$ws.Range("J112").Value = 2415.7144
$ws.Range("L112").Value = 7247.1432
$ws.Range("N112").Value = -9463.143199999999
# Row 121
$ws.Range("H121").Value = 620.13336
$ws.Range("I121").Value = 2350
$ws.Range("J121").Value = 354
$ws.Range("K121").Value = 7050
$ws.Range("L121").Value = 1062
$ws.Range("M121").Value = -5303
$ws.Range("N121").Value = -4556
# Row 141
$ws.Range("H141").Value = 3799
$ws.Range("I141").Value = 2190.7693
$ws.Range("K141").Value = 6572.3079
$ws.Range("M141").Value = -1392.3079

$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 13666.667
$ws.Range("J43").Value = 13666.667
$ws.Range("L43").Value = 13666.667
$ws.Range("N43").Value = -14292.667
# Row 45
$ws.Range("H45").Value = 772063.4399999999
$ws.Range("I45").Value = 1001281.1
$ws.Range("K45").Value = 1001281.1
$ws.Range("M45").Value = -1000904.1
# Row 102
$ws.Range("H102").Value = 2210
$ws.Range("I102").Value = 2615.5557
$ws.Range("J102").Value = 1688.5714
$ws.Range("K102").Value = 2615.5557
$ws.Range("L102").Value = 1688.5714
$ws.Range("M102").Value = -993.5556999999999
$ws.Range("N102").Value = -4932.5714

$ws = $wb.Worksheets.Item("BSM")
# Row 60
$ws.Range("H60").Value = 22490
$ws.Range("J60").Value = 22490
$ws.Range("L60").Value = 22490
$ws.Range("N60").Value = -23688
# Row 99
$ws.Range("H99").Value = 728.3570999999999
$ws.Range("I99").Value = 707.46155
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 707.46155
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 790.53845
$ws.Range("N99").Value = -3996
# Row 103
$ws.Range("H103").Value = 19500
$ws.Range("J103").Value = 19500
$ws.Range("L103").Value = 19500
$ws.Range("N103").Value = -21844
# Row 105
$ws.Range("H105").Value = 1516.6333
$ws.Range("I105").Value = 1505.7916
$ws.Range("J105").Value = 1560
$ws.Range("K105").Value = 1505.7916
$ws.Range("L105").Value = 1560
$ws.Range("M105").Value = 241.2084
$ws.Range("N105").Value = -5054
# Row 111
$ws.Range("H111").Value = 24702
$ws.Range("J111").Value = 24702
$ws.Range("L111").Value = 24702
$ws.Range("N111").Value = -32882
# Row 134
$ws.Range("H134").Value = 13575748
$ws.Range("I134").Value = 14286754
$ws.Range("J134").Value = 7354444.5
$ws.Range("K134").Value = 42860262
$ws.Range("L134").Value = 22063333.5
$ws.Range("M134").Value = -42857727
$ws.Range("N134").Value = -22068403.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 15386385
$ws.Range("I31").Value = 23810556
$ws.Range("J31").Value = 3115.0435
$ws.Range("K31").Value = 23810556
$ws.Range("L31").Value = 3115.0435
$ws.Range("M31").Value = -23810261
$ws.Range("N31").Value = -3705.0435
# Row 34
$ws.Range("H34").Value = 15386385
$ws.Range("I34").Value = 23810556
$ws.Range("J34").Value = 3115.0435
$ws.Range("K34").Value = 23810556
$ws.Range("L34").Value = 3115.0435
$ws.Range("M34").Value = -23810354
$ws.Range("N34").Value = -3519.0435
# Row 86
$ws.Range("H86").Value = 9411.806
$ws.Range("I86").Value = 11365.8
$ws.Range("J86").Value = 4970.909
$ws.Range("K86").Value = 11365.8
$ws.Range("L86").Value = 4970.909
$ws.Range("M86").Value = -10242.8
$ws.Range("N86").Value = -7216.909
# Row 89
$ws.Range("H89").Value = 9411.806
$ws.Range("I89").Value = 11365.8
$ws.Range("J89").Value = 4970.909
$ws.Range("K89").Value = 56829
$ws.Range("L89").Value = 24854.545
$ws.Range("M89").Value = -51213
$ws.Range("N89").Value = -36086.545
# Row 107
$ws.Range("H107").Value = 577.3611
$ws.Range("I107").Value = 309.5238
$ws.Range("J107").Value = 952.3333
$ws.Range("K107").Value = 309.5238
$ws.Range("L107").Value = 952.3333
$ws.Range("M107").Value = 1610.4762
$ws.Range("N107").Value = -4792.3333
# Row 122
$ws.Range("H122").Value = 7259.4
$ws.Range("I122").Value = 13701.714
$ws.Range("J122").Value = 1622.375
$ws.Range("K122").Value = 41105.142
$ws.Range("L122").Value = 4867.125
$ws.Range("M122").Value = -38655.142
$ws.Range("N122").Value = -9767.125

$ws = $wb.Worksheets.Item("CUL")
# Row 133
$ws.Range("H133").Value = 2725.8333
$ws.Range("I133").Value = 2725.8333
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 8177.499899999999
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -3117.499899999999
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 16674899
$ws.Range("I97").Value = 1178.3636
$ws.Range("J97").Value = 62527630
$ws.Range("K97").Value = 1178.3636
$ws.Range("L97").Value = 62527630
$ws.Range("M97").Value = -682.3635999999999
$ws.Range("N97").Value = -62528622
# Row 111
$ws.Range("H111").Value = 20293
$ws.Range("J111").Value = 20293
$ws.Range("L111").Value = 20293
$ws.Range("N111").Value = -26427

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4169866.5
$ws.Range("I61").Value = 20833332
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 20833332
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -20833130
$ws.Range("N61").Value = -4404
# Row 100
$ws.Range("H100").Value = 1564.7059
$ws.Range("J100").Value = 2043.75
$ws.Range("L100").Value = 2043.75
$ws.Range("N100").Value = -3125.75
# Row 113
$ws.Range("H113").Value = 4169866.5
$ws.Range("I113").Value = 20833332
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 20833332
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -20831162
$ws.Range("N113").Value = -8340

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 509115.88
$ws.Range("I132").Value = 605328.5
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 1815985.5
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -1813455.5
$ws.Range("N132").Value = -17058.5

Write-Host "Applied Ridill_Profits price updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets."
